# Fruta / hortaliza, semanal
# Insert two new weekly-report rows at the top of this date-ordered block
# (rows 573-662), pushing the existing data down by two rows, and fill
# the two newly inserted rows with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 573 (existing row 573 and
# everything below it shifts down to 575.. and beyond).
$ws.Rows.Item(573).Insert()
$ws.Rows.Item(573).Insert()

# New row 573: Larga vida / Primera
$ws.Range("A573").Value = 3
$ws.Range("B573").Value = "Femacal de La Calera"
$ws.Range("C573").Value = "Coquimbo"
$ws.Range("D573").Value = 44474
$ws.Range("E573").Value = 5
$ws.Range("F573").Value = 100112020
$ws.Range("G573").Value = "Tomate"
$ws.Range("H573").Value = "Larga vida"
$ws.Range("I573").Value = "Primera"
$ws.Range("J573").Value = 290
$ws.Range("K573").Value = 17000
$ws.Range("L573").Value = 18000
$ws.Range("M573").Value = 17517
$ws.Range("N573").Value = "`$/bandeja 18 kilos"
$ws.Range("O573").Value = "Región de Arica y Parinacota"
$ws.Range("P573").Value = 973
$ws.Range("Q573").Value = 18
$ws.Range("R573").Value = "Hortaliza"

# New row 574: Larga vida / Segunda
$ws.Range("A574").Value = 3
$ws.Range("B574").Value = "Femacal de La Calera"
$ws.Range("C574").Value = "Coquimbo"
$ws.Range("D574").Value = 44474
$ws.Range("E574").Value = 5
$ws.Range("F574").Value = 100112020
$ws.Range("G574").Value = "Tomate"
$ws.Range("H574").Value = "Larga vida"
$ws.Range("I574").Value = "Segunda"
$ws.Range("J574").Value = 320
$ws.Range("K574").Value = 14000
$ws.Range("L574").Value = 15000
$ws.Range("M574").Value = 14500
$ws.Range("N574").Value = "`$/bandeja 18 kilos"
$ws.Range("O574").Value = "Región de Arica y Parinacota"
$ws.Range("P574").Value = 806
$ws.Range("Q574").Value = 18
$ws.Range("R574").Value = "Hortaliza"

# Apply the date number format (style used by column D elsewhere) to the
# two new date cells.
$ws.Range("D573:D574").NumberFormat = "YYYY-MM-DD HH:MM:SS"
